# Auto-generated Excel COM-interop script
# Applies refreshed market-data values to the Hades Profits workbook (per sheet: ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR)
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 34: Sophomore Slump | Goatskin Grimoire
$ws.Range("H34").Value = 3000
$ws.Range("I34").Value = 3000
$ws.Range("K34").Value = 3000
$ws.Range("M34").Value = -2797
# Row 36: You Put Your Left Hand In | Engraved Goatskin Grimoire
$ws.Range("H36").Value = 3000
$ws.Range("I36").Value = 3000
$ws.Range("K36").Value = 3000
$ws.Range("M36").Value = -2285
# Row 40: Stuck in the Moment | Horn Glue
$ws.Range("H40").Value = 3171.4285
$ws.Range("I40").Value = 2066.6667
$ws.Range("J40").Value = 4000
$ws.Range("K40").Value = 2066.6667
$ws.Range("L40").Value = 4000
$ws.Range("M40").Value = -1891.6667
$ws.Range("N40").Value = -4350
# Row 112: Making Ends Meet | Superior Spiritbond Potion
$ws.Range("H112").Value = 27779162
$ws.Range("I112").Value = 250000420
$ws.Range("J112").Value = 3087911
$ws.Range("K112").Value = 750001260
$ws.Range("L112").Value = 9263733
$ws.Range("M112").Value = -750000152
$ws.Range("N112").Value = -9265949

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 2: Ain't Got No Ingots | Bronze Ingot
$ws.Range("H2").Value = 2108.7036
$ws.Range("I2").Value = 1941.6
$ws.Range("K2").Value = 1941.6
$ws.Range("M2").Value = -1828.6
# Row 32: Ingot We Trust | Steel Ingot
$ws.Range("H32").Value = 14416.424
$ws.Range("I32").Value = 15148.91
$ws.Range("J32").Value = 10335.429
$ws.Range("K32").Value = 15148.91
$ws.Range("L32").Value = 10335.429
$ws.Range("M32").Value = -14861.91
$ws.Range("N32").Value = -10909.429
# Row 116: No Scope | Titanbronze Ingot
$ws.Range("H116").Value = 2108.7036
$ws.Range("I116").Value = 1941.6
$ws.Range("K116").Value = 1941.6
$ws.Range("M116").Value = 352.4000000000001
# Row 122: Haste for High Durium | High Durium Nugget
$ws.Range("H122").Value = 5557821.5
$ws.Range("I122").Value = 2221.0667
$ws.Range("J122").Value = 22224622
$ws.Range("K122").Value = 6663.2001
$ws.Range("L122").Value = 66673866
$ws.Range("M122").Value = -4213.2001
$ws.Range("N122").Value = -66678766

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 3: Hells Bells | Bronze Ingot
$ws.Range("H3").Value = 2108.7036
$ws.Range("I3").Value = 1941.6
$ws.Range("K3").Value = 1941.6
$ws.Range("M3").Value = -1827.6

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31: Wall Not Found | Walnut Lumber
$ws.Range("H31").Value = 2647.325
$ws.Range("I31").Value = 1418.2963
$ws.Range("J31").Value = 5199.923
$ws.Range("K31").Value = 1418.2963
$ws.Range("L31").Value = 5199.923
$ws.Range("M31").Value = -1123.2963
$ws.Range("N31").Value = -5789.923
# Row 34: Armoires of the Rich and Famous | Walnut Lumber
$ws.Range("H34").Value = 2647.325
$ws.Range("I34").Value = 1418.2963
$ws.Range("J34").Value = 5199.923
$ws.Range("K34").Value = 1418.2963
$ws.Range("L34").Value = 5199.923
$ws.Range("M34").Value = -1216.2963
$ws.Range("N34").Value = -5603.923
# Row 108: Just Starting Out | White Oak Fishing Rod
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").ClearContents()
# Row 134: Wood You Be Quiet | Ceiba Lumber
$ws.Range("H134").Value = 67086.12
$ws.Range("I134").Value = 4189.5386
$ws.Range("J134").Value = 271500
$ws.Range("K134").Value = 12568.6158
$ws.Range("L134").Value = 814500
$ws.Range("M134").Value = -10033.6158
$ws.Range("N134").Value = -819570

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 16: Go Ahead and Dig In | Mole Loaf
$ws.Range("H16").Value = 1911.5385
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 1911.5385
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 5734.6155
$ws.Range("N16").Value = -6080.6155
$ws.Range("M16").ClearContents()
# Row 68: Such a Butter Face | Fermented Butter
$ws.Range("H68").Value = 409163.97
$ws.Range("I68").Value = 964.9091
$ws.Range("J68").Value = 1251074.5
$ws.Range("K68").Value = 2894.7273
$ws.Range("L68").Value = 3753223.5
$ws.Range("M68").Value = -2083.7273
$ws.Range("N68").Value = -3754845.5
# Row 71: No Margarine of Error (L) | Fermented Butter
$ws.Range("H71").Value = 409163.97
$ws.Range("I71").Value = 964.9091
$ws.Range("J71").Value = 1251074.5
$ws.Range("K71").Value = 8684.1819
$ws.Range("L71").Value = 11259670.5
$ws.Range("M71").Value = -4628.1819
$ws.Range("N71").Value = -11267782.5
# Row 103: West Meats East | Nomad Meat Pie
$ws.Range("H103").Value = 2901.9
$ws.Range("I103").Value = 506.25
$ws.Range("J103").Value = 3773.0454
$ws.Range("K103").Value = 1518.75
$ws.Range("L103").Value = 11319.1362
$ws.Range("M103").Value = -639.75
$ws.Range("N103").Value = -13077.1362
# Row 107: Slippery Service | Frantoio Oil
$ws.Range("H107").Value = 642.2083
$ws.Range("I107").Value = 1077.9333
$ws.Range("J107").Value = 444.15152
$ws.Range("K107").Value = 3233.7999
$ws.Range("L107").Value = 1332.45456
$ws.Range("M107").Value = -1313.7999
$ws.Range("N107").Value = -5172.45456
# Row 131: The Mountain Steeped | Tsai tou Vounou
$ws.Range("H131").Value = 967.6824
$ws.Range("J131").Value = 1031.9333
$ws.Range("L131").Value = 3095.7999
$ws.Range("N131").Value = -13175.7999

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 70: Sky Is the Limit | Mythrite Ingot
$ws.Range("H70").Value = 34977.44
$ws.Range("I70").Value = 52968.24
$ws.Range("K70").Value = 52968.24
$ws.Range("M70").Value = -52698.24
# Row 73: Hulls of Broken Dreams (L) | Mythrite Ingot
$ws.Range("H73").Value = 34977.44
$ws.Range("I73").Value = 52968.24
$ws.Range("K73").Value = 52968.24
$ws.Range("M73").Value = -52032.24
# Row 101: Best-laid Planispheres | Dual-plated Durium Planisphere
$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents()
# Row 122: Awarding Academic Excellence | Ametrine
$ws.Range("H122").Value = 2836.0625
$ws.Range("I122").Value = 2016.9
$ws.Range("J122").Value = 4201.3335
$ws.Range("K122").Value = 6050.700000000001
$ws.Range("L122").Value = 12604.0005
$ws.Range("M122").Value = -3600.700000000001
$ws.Range("N122").Value = -17504.0005

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 7: Tan Before the Ban | Leather
$ws.Range("H7").Value = 2400
$ws.Range("I7").Value = 2355.5557
$ws.Range("J7").Value = 2800
$ws.Range("K7").Value = 2355.5557
$ws.Range("L7").Value = 2800
$ws.Range("M7").Value = -2243.5557
$ws.Range("N7").Value = -3024
# Row 22: Skin off Their Backs | Aldgoat Leather
$ws.Range("H22").Value = 1267.1111
$ws.Range("I22").Value = 400
$ws.Range("J22").Value = 1375.5
$ws.Range("K22").Value = 400
$ws.Range("L22").Value = 1375.5
$ws.Range("M22").Value = -105
$ws.Range("N22").Value = -1965.5
# Row 27: Fire and Hide | Aldgoat Leather
$ws.Range("H27").Value = 1267.1111
$ws.Range("I27").Value = 400
$ws.Range("J27").Value = 1375.5
$ws.Range("K27").Value = 400
$ws.Range("L27").Value = 1375.5
$ws.Range("M27").Value = -293
$ws.Range("N27").Value = -1589.5
# Row 55: It's Not a Job, It's a Calling | Peiste Leather
$ws.Range("H55").Value = 354.3846
$ws.Range("J55").Value = 398.6
$ws.Range("L55").Value = 398.6
$ws.Range("N55").Value = -744.6
# Row 126: Battered Books | Saiga Leather
$ws.Range("H126").Value = 2400
$ws.Range("I126").Value = 2355.5557
$ws.Range("J126").Value = 2800
$ws.Range("K126").Value = 7066.6671
$ws.Range("L126").Value = 8400
$ws.Range("M126").Value = -4596.6671
$ws.Range("N126").Value = -13340
# Row 132: Tenets of Tanning | Silver Lobo Leather
$ws.Range("H132").Value = 46896.293
$ws.Range("I132").Value = 22479.04
$ws.Range("J132").Value = 114722
$ws.Range("K132").Value = 67437.12
$ws.Range("L132").Value = 344166
$ws.Range("M132").Value = -64907.12
$ws.Range("N132").Value = -349226

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 14: Hat in Hand | Straw Hat
$ws.Range("H14").Value = 9483.223
$ws.Range("J14").Value = 6918.75
$ws.Range("L14").Value = 6918.75
$ws.Range("N14").Value = -7254.75
# Row 132: Comfy Cabins | Snow Cotton Cloth
$ws.Range("H132").Value = 85055.5
$ws.Range("I132").Value = 60078.47
$ws.Range("J132").Value = 145714
$ws.Range("K132").Value = 180235.41
$ws.Range("L132").Value = 437142
$ws.Range("M132").Value = -177705.41
$ws.Range("N132").Value = -442202
# Row 133: Begin with the Basics | Snow Cotton Jacket
$ws.Range("H133").Value = 37900
$ws.Range("J133").Value = 37900
$ws.Range("L133").Value = 37900
$ws.Range("N133").Value = -48020
